# QuickQuoteTestData.xlsx update
# Refresh stale shipment test data (pickup date, order ids, clone order id,
# tracking numbers and waybill numbers) on the "Input" sheet for rows 3 and 10
# so that QA / staging shipment tests stop failing against expired references.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Helper: write a value into a cell as literal text, avoiding Excel's
# automatic number/date inference (so e.g. "02-12-2022" or "51530846" are
# kept as text shared strings instead of being converted into a date serial
# or a numeric value), while preserving the cell's existing style.
function Set-TextValue {
    param(
        $Cell,
        [string]$Value
    )
    $Cell.Formula = '=TEXT("' + $Value + '","@")'
    $Cell.Copy() | Out-Null
    $Cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 3
Set-TextValue $ws.Range("B3") "02-12-2022"
Set-TextValue $ws.Range("V3") "51530846"
Set-TextValue $ws.Range("Y3") "FCT942085809463361536"
Set-TextValue $ws.Range("Z3") "CEV1003987"

# Row 10
Set-TextValue $ws.Range("B10") "02-12-2022"
Set-TextValue $ws.Range("V10") "51530843"
Set-TextValue $ws.Range("W10") "51530845"
Set-TextValue $ws.Range("Y10") "1Z44R7R60394589476"
Set-TextValue $ws.Range("Z10") "FCUPSG1012088"

$excel.CutCopyMode = $false
